# Auto-generated Excel COM-interop script
# Applies scraped market-price / profit updates to the Garuda_Profits-style workbook.
# Each block targets one Leve row on one job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR),
# updating currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) to freshly
# scraped values. Some rows gain or lose a LeveProfitHQ (column N) cell depending on
# whether an HQ price now exists.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 8: On the Drip
$ws.Range("H8").Value = 47.25
$ws.Range("I8").Value = 47.25
$ws.Range("K8").Value = 141.75
$ws.Range("M8").Value = -2.75

# Row 21: Book and a Hard Place
$ws.Range("H21").Value = 7874
$ws.Range("I21").Value = 1345
$ws.Range("K21").Value = 1345
$ws.Range("M21").Value = -877

# Row 23: There's Something about Bury
$ws.Range("H23").Value = 7874
$ws.Range("I23").Value = 1345
$ws.Range("K23").Value = 1345
$ws.Range("M23").Value = -1111

# Row 29: Dripping with Venom
$ws.Range("H29").Value = 2600.3333
$ws.Range("I29").Value = 866.6667
$ws.Range("J29").Value = 3467.1667
$ws.Range("K29").Value = 2600.0001
$ws.Range("L29").Value = 10401.5001
$ws.Range("M29").Value = -2319.0001
$ws.Range("N29").Value = -10963.5001

# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 80.2
$ws.Range("I38").Value = 80.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 240.6
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 131.4
$ws.Range("N38").ClearContents()

# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 6203.84
$ws.Range("I58").Value = 261.23077
$ws.Range("K58").Value = 783.69231
$ws.Range("M58").Value = -633.69231

# Row 63: Summoning for Dummies
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66: Summoning the Courage to Be Different (L)
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 138.13333
$ws.Range("I107").Value = 138.13333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 138.13333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1781.86667
$ws.Range("N107").ClearContents()

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1579.3334
$ws.Range("J112").Value = 1617.64
$ws.Range("L112").Value = 4852.92
$ws.Range("N112").Value = -7068.92

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1779.641
$ws.Range("I137").Value = 1399.7693
$ws.Range("J137").Value = 2539.3845
$ws.Range("K137").Value = 4199.3079
$ws.Range("L137").Value = 7618.1535
$ws.Range("M137").Value = -1649.3079
$ws.Range("N137").Value = -12718.1535

$ws = $wb.Worksheets.Item("ARM")

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1058.6666
$ws.Range("I61").Value = 671.5833
$ws.Range("K61").Value = 671.5833
$ws.Range("M61").Value = -459.5833

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 918
$ws.Range("I74").Value = 822
$ws.Range("K74").Value = 822
$ws.Range("M74").Value = 52

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 918
$ws.Range("I77").Value = 822
$ws.Range("K77").Value = 4110
$ws.Range("M77").Value = 258

# Row 94: Setting the Stage
$ws.Range("H94").Value = 34830
$ws.Range("J94").Value = 34830
$ws.Range("L94").Value = 34830
$ws.Range("N94").Value = -36632

# Row 96: The Gauntlet Is Cast
$ws.Range("H96").Value = 17984.143
$ws.Range("J96").Value = 17984.143
$ws.Range("L96").Value = 17984.143
$ws.Range("N96").Value = -23476.143

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2867.7
$ws.Range("I132").Value = 2814.1865
$ws.Range("J132").Value = 3154.7273
$ws.Range("K132").Value = 8442.559499999999
$ws.Range("L132").Value = 9464.1819
$ws.Range("M132").Value = -5912.559499999999
$ws.Range("N132").Value = -14524.1819

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1058.6666
$ws.Range("I136").Value = 671.5833
$ws.Range("K136").Value = 2014.7499
$ws.Range("M136").Value = 535.2501

$ws = $wb.Worksheets.Item("BSM")

# Row 53: Kitchen Casualties
$ws.Range("H53").Value = 45000
$ws.Range("J53").Value = 45000
$ws.Range("L53").Value = 45000
$ws.Range("N53").Value = -46148

# Row 100: And My Axe
$ws.Range("H100").Value = 32630.75
$ws.Range("J100").Value = 32630.75
$ws.Range("L100").Value = 32630.75
$ws.Range("N100").Value = -34794.75

# Row 133: Paring Is Caring
$ws.Range("H133").Value = 47000
$ws.Range("J133").Value = 47000
$ws.Range("L133").Value = 47000
$ws.Range("N133").Value = -57120

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 19871.455
$ws.Range("I134").Value = 24166.727
$ws.Range("J134").Value = 2690.3635
$ws.Range("K134").Value = 72500.181
$ws.Range("L134").Value = 8071.0905
$ws.Range("M134").Value = -69965.181
$ws.Range("N134").Value = -13141.0905

$ws = $wb.Worksheets.Item("CRP")

# Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 76.166664
$ws.Range("I7").Value = 68.666664
$ws.Range("J7").Value = 83.666664
$ws.Range("K7").Value = 68.666664
$ws.Range("L7").Value = 83.666664
$ws.Range("M7").Value = 44.333336
$ws.Range("N7").Value = -309.666664

# Row 31: Wall Not Found
$ws.Range("H31").Value = 4275524.5
$ws.Range("I31").Value = 1367.9333
$ws.Range("J31").Value = 18522712
$ws.Range("K31").Value = 1367.9333
$ws.Range("L31").Value = 18522712
$ws.Range("M31").Value = -1072.9333
$ws.Range("N31").Value = -18523302

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 4275524.5
$ws.Range("I34").Value = 1367.9333
$ws.Range("J34").Value = 18522712
$ws.Range("K34").Value = 1367.9333
$ws.Range("L34").Value = 18522712
$ws.Range("M34").Value = -1165.9333
$ws.Range("N34").Value = -18523116

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 9524731
$ws.Range("I58").Value = 966.96155
$ws.Range("J58").Value = 37037828
$ws.Range("K58").Value = 966.96155
$ws.Range("L58").Value = 37037828
$ws.Range("M58").Value = -763.96155
$ws.Range("N58").Value = -37038234

# Row 102: The Ear Is the Way to the Heart
$ws.Range("H102").Value = 41875
$ws.Range("J102").Value = 41875
$ws.Range("L102").Value = 41875
$ws.Range("N102").Value = -46743

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2091.225
$ws.Range("I132").Value = 2042.6129
$ws.Range("J132").Value = 2258.6667
$ws.Range("K132").Value = 6127.8387
$ws.Range("L132").Value = 6776.000100000001
$ws.Range("M132").Value = -3597.8387
$ws.Range("N132").Value = -11836.0001

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1206.641
$ws.Range("I134").Value = 1178.2059
$ws.Range("K134").Value = 3534.6177
$ws.Range("M134").Value = -999.6176999999998

# Row 136: Turali Quality
$ws.Range("H136").Value = 9524731
$ws.Range("I136").Value = 966.96155
$ws.Range("J136").Value = 37037828
$ws.Range("K136").Value = 2900.88465
$ws.Range("L136").Value = 111113484
$ws.Range("M136").Value = -350.88465
$ws.Range("N136").Value = -111118584

$ws = $wb.Worksheets.Item("CUL")

# Row 68: Such a Butter Face
$ws.Range("H68").Value = 293
$ws.Range("I68").Value = 310.16666
$ws.Range("J68").Value = 190
$ws.Range("K68").Value = 930.4999799999999
$ws.Range("L68").Value = 570
$ws.Range("M68").Value = -119.4999799999999
$ws.Range("N68").Value = -2192

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 293
$ws.Range("I71").Value = 310.16666
$ws.Range("J71").Value = 190
$ws.Range("K71").Value = 2791.49994
$ws.Range("L71").Value = 1710
$ws.Range("M71").Value = 1264.50006
$ws.Range("N71").Value = -9822

# Row 130: Blast from the Pasta
$ws.Range("H130").Value = 1226.25
$ws.Range("I130").Value = 901.4286
$ws.Range("J130").Value = 3500
$ws.Range("K130").Value = 2704.2858
$ws.Range("L130").Value = 10500
$ws.Range("M130").Value = 2315.7142
$ws.Range("N130").Value = -20540

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 820181.9399999999
$ws.Range("J131").Value = 1905788.5
$ws.Range("L131").Value = 5717365.5
$ws.Range("N131").Value = -5727445.5

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 30914730
$ws.Range("I137").Value = 27779578
$ws.Range("J137").Value = 33422852
$ws.Range("K137").Value = 83338734
$ws.Range("L137").Value = 100268556
$ws.Range("M137").Value = -83333634
$ws.Range("N137").Value = -100278756

# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 2047.9565
$ws.Range("I140").Value = 1431.6666
$ws.Range("J140").Value = 4266.6
$ws.Range("K140").Value = 4294.9998
$ws.Range("L140").Value = 12799.8
$ws.Range("M140").Value = 885.0002000000004
$ws.Range("N140").Value = -23159.8

$ws = $wb.Worksheets.Item("LTW")

# Row 18: Simply the Best
$ws.Range("H18").Value = 11667.667
$ws.Range("J18").Value = 11667.667
$ws.Range("L18").Value = 11667.667
$ws.Range("N18").Value = -12011.667

# Row 104: Brace Yourselves
$ws.Range("H104").Value = 29978.5
$ws.Range("J104").Value = 29978.5
$ws.Range("L104").Value = 29978.5
$ws.Range("N104").Value = -36966.5

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 6281.143
$ws.Range("I132").Value = 8404.25
$ws.Range("J132").Value = 2034.9286
$ws.Range("K132").Value = 25212.75
$ws.Range("L132").Value = 6104.7858
$ws.Range("M132").Value = -22682.75
$ws.Range("N132").Value = -11164.7858

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4726.1934
$ws.Range("I136").Value = 5554.864
$ws.Range("J136").Value = 2700.5557
$ws.Range("K136").Value = 16664.592
$ws.Range("L136").Value = 8101.6671
$ws.Range("M136").Value = -14114.592
$ws.Range("N136").Value = -13201.6671

$ws = $wb.Worksheets.Item("WVR")

# Row 12: This Is Why You Can't Have Nice Things
$ws.Range("H12").Value = 80007
$ws.Range("J12").Value = 80007
$ws.Range("L12").Value = 80007
$ws.Range("N12").Value = -80291

# Row 101: Who War It Better
$ws.Range("H101").Value = 9500
$ws.Range("J101").Value = 9500
$ws.Range("L101").Value = 9500
$ws.Range("N101").Value = -15990

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 858.0923
$ws.Range("I132").Value = 819.0714
$ws.Range("J132").Value = 1100.8889
$ws.Range("K132").Value = 2457.2142
$ws.Range("L132").Value = 3302.6667
$ws.Range("M132").Value = 72.78579999999965
$ws.Range("N132").Value = -8362.6667

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 3027.5178
$ws.Range("I136").Value = 3336.8696
$ws.Range("J136").Value = 1604.5
$ws.Range("K136").Value = 10010.6088
$ws.Range("L136").Value = 4813.5
$ws.Range("M136").Value = -7460.6088
$ws.Range("N136").Value = -9913.5

# Row 138: Halfgloves, Full Effort
$ws.Range("H138").Value = 51389.855
$ws.Range("J138").Value = 51389.855
$ws.Range("L138").Value = 51389.855
$ws.Range("N138").Value = -61669.855
